$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.776.58'
$ws.Range('E2').Value = '  -7.17%  '
$ws.Range('D3').Value = '2.543.78'
$ws.Range('E3').Value = '  -2.20%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '296.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '91.51'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.30%  '
$ws.Range('E7').Value = '  -4.34%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -6.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.63'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -8.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0804'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.65'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -6.01%  '
$ws.Range('D13').Value = '2.936.24'
$ws.Range('E13').Value = '  -2.10%  '
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('D15').Value = '2.540.29'
$ws.Range('E15').Value = '  -2.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.861'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -6.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.09'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.32%  '
$ws.Range('D18').Value = '42.805.23'
$ws.Range('E18').Value = '  -7.38%  '
$ws.Range('E19').Value = '  -4.38%  '
$ws.Range('E20').Value = '  -1.40%  '
$ws.Range('E21').Value = '  -2.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.37'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '259.58'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -12.07%  '
$ws.Range('E24').Value = '  -6.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '29.46'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.35%  '
$ws.Range('E26').Value = '  -7.57%  '
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('E28').Value = '  -7.48%  '
$ws.Range('E29').Value = '  -4.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.04'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.88'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '150.42'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.15'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.37'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.73%  '
$ws.Range('E35').Value = '  -3.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0790'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.84%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.113'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.10'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +11.51%  '
$ws.Range('E39').Value = '  -3.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.10'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.41'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.27%  '
$ws.Range('E42').Value = '  -6.87%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.80'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.72%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.073.23'
$ws.Range('E44').Value = '  -1.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '84.89'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -13.53%  '
$ws.Range('E47').Value = '  +2.86%  '
$ws.Range('E48').Value = '  -2.02%  '
$ws.Range('D49').Value = '2.792.00'
$ws.Range('E49').Value = '  -2.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.66'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -9.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '103.61'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.47%  '
